$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 22 <- former row 26 data
$ws.Range("B22").Value = 6880327
$ws.Range("E22").Value = "Dynamo Dresden"
$ws.Range("F22").Value = "Waldhof Mannheim"
$ws.Range("G22").Value = 2
$ws.Range("H22").Value = 1
$ws.Range("I22").Value = 2
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = "H"
$ws.Range("L22").Value = 1.65
$ws.Range("M22").Value = 3.8
$ws.Range("N22").Value = 4.333
$ws.Range("O22").Value = 1.65
$ws.Range("P22").Value = 4
$ws.Range("Q22").Value = 4
$ws.Range("R22").Value = -0.75
$ws.Range("S22").Value = 1.8
$ws.Range("T22").Value = 2
$ws.Range("U22").Value = 3
$ws.Range("V22").Value = 1.85
$ws.Range("W22").Value = 1.95
$ws.Range("X22").Value = 0.6499999999999999
$ws.Range("Y22").Value = -1
$ws.Range("Z22").Value = -1
$ws.Range("AA22").Value = 0.4
$ws.Range("AB22").Value = -0.5
$ws.Range("AC22").Value = 0
$ws.Range("AD22").Value = 0

# Row 23 <- former row 22 data
$ws.Range("B23").Value = 6880508
$ws.Range("E23").Value = "Hallescher FC"
$ws.Range("F23").Value = "MSV Duisburg"
$ws.Range("G23").Value = 1
$ws.Range("H23").Value = 1
$ws.Range("I23").Value = 1
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = "D"
$ws.Range("L23").Value = 2.375
$ws.Range("M23").Value = 3.5
$ws.Range("N23").Value = 2.6
$ws.Range("O23").Value = 2.15
$ws.Range("P23").Value = 3.6
$ws.Range("Q23").Value = 2.8
$ws.Range("R23").Value = -0.25
$ws.Range("S23").Value = 1.95
$ws.Range("T23").Value = 1.85
$ws.Range("U23").Value = 2.75
$ws.Range("V23").Value = 1.8
$ws.Range("W23").Value = 2
$ws.Range("X23").Value = -1
$ws.Range("Y23").Value = 2.6
$ws.Range("Z23").Value = -1
$ws.Range("AA23").Value = -0.5
$ws.Range("AB23").Value = 0.425
$ws.Range("AC23").Value = -1
$ws.Range("AD23").Value = 1

# Row 24 <- former row 23 data
$ws.Range("B24").Value = 6881313
$ws.Range("E24").Value = "SSV Ulm 1846"
$ws.Range("F24").Value = "Arminia Bielefeld"
$ws.Range("G24").Value = 1
$ws.Range("H24").Value = 0
$ws.Range("I24").Value = 1
$ws.Range("J24").Value = 0
$ws.Range("K24").Value = "H"
$ws.Range("L24").Value = 3.2
$ws.Range("M24").Value = 3.4
$ws.Range("N24").Value = 2.05
$ws.Range("O24").Value = 2.8
$ws.Range("P24").Value = 3.5
$ws.Range("Q24").Value = 2.25
$ws.Range("R24").Value = 0.25
$ws.Range("S24").Value = 1.825
$ws.Range("T24").Value = 1.975
$ws.Range("U24").Value = 3
$ws.Range("V24").Value = 1.95
$ws.Range("W24").Value = 1.85
$ws.Range("X24").Value = 1.8
$ws.Range("Y24").Value = -1
$ws.Range("Z24").Value = -1
$ws.Range("AA24").Value = 0.825
$ws.Range("AB24").Value = -1
$ws.Range("AC24").Value = -1
$ws.Range("AD24").Value = 0.8500000000000001

# Row 25 <- former row 24 data
$ws.Range("B25").Value = 6880353
$ws.Range("E25").Value = "SC Preussen Munster"
$ws.Range("F25").Value = "FC Ingolstadt"
$ws.Range("G25").Value = 3
$ws.Range("H25").Value = 1
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 1
$ws.Range("K25").Value = "H"
$ws.Range("L25").Value = 2.3
$ws.Range("M25").Value = 3.4
$ws.Range("N25").Value = 2.8
$ws.Range("O25").Value = 2.25
$ws.Range("P25").Value = 3.2
$ws.Range("Q25").Value = 2.8
$ws.Range("R25").Value = -0.25
$ws.Range("S25").Value = 2.025
$ws.Range("T25").Value = 1.775
$ws.Range("U25").Value = 2.75
$ws.Range("V25").Value = 1.9
$ws.Range("W25").Value = 1.9
$ws.Range("X25").Value = 1.25
$ws.Range("Y25").Value = -1
$ws.Range("Z25").Value = -1
$ws.Range("AA25").Value = 1.025
$ws.Range("AB25").Value = -1
$ws.Range("AC25").Value = 0.8999999999999999
$ws.Range("AD25").Value = -1

# Row 26 <- former row 25 data
$ws.Range("B26").Value = 6881596
$ws.Range("E26").Value = "1860 Munich"
$ws.Range("F26").Value = "Vfb Lubeck"
$ws.Range("G26").Value = 1
$ws.Range("H26").Value = 2
$ws.Range("I26").Value = 1
$ws.Range("J26").Value = 2
$ws.Range("K26").Value = "A"
$ws.Range("L26").Value = 1.615
$ws.Range("M26").Value = 3.75
$ws.Range("N26").Value = 4.75
$ws.Range("O26").Value = 1.615
$ws.Range("P26").Value = 4
$ws.Range("Q26").Value = 4.5
$ws.Range("R26").Value = -0.75
$ws.Range("S26").Value = 1.775
$ws.Range("T26").Value = 2.025
$ws.Range("U26").Value = 3
$ws.Range("V26").Value = 1.95
$ws.Range("W26").Value = 1.85
$ws.Range("X26").Value = -1
$ws.Range("Y26").Value = -1
$ws.Range("Z26").Value = 3.5
$ws.Range("AA26").Value = -1
$ws.Range("AB26").Value = 1.025
$ws.Range("AC26").Value = 0
$ws.Range("AD26").Value = 0

# Row 29 <- former row 30 data
$ws.Range("B29").Value = 6881594
$ws.Range("E29").Value = "Saarbrucken"
$ws.Range("F29").Value = "Verl"
$ws.Range("G29").Value = 4
$ws.Range("H29").Value = 3
$ws.Range("I29").Value = 1
$ws.Range("J29").Value = 2
$ws.Range("K29").Value = "H"
$ws.Range("L29").Value = 1.6
$ws.Range("M29").Value = 3.8
$ws.Range("N29").Value = 4.75
$ws.Range("O29").Value = 1.75
$ws.Range("P29").Value = 3.8
$ws.Range("Q29").Value = 3.8
$ws.Range("R29").Value = -0.75
$ws.Range("S29").Value = 2
$ws.Range("T29").Value = 1.8
$ws.Range("U29").Value = 3
$ws.Range("V29").Value = 1.85
$ws.Range("W29").Value = 1.95
$ws.Range("X29").Value = 0.75
$ws.Range("Y29").Value = -1
$ws.Range("Z29").Value = -1
$ws.Range("AA29").Value = 0.5
$ws.Range("AB29").Value = -0.5
$ws.Range("AC29").Value = 0.8500000000000001
$ws.Range("AD29").Value = -1

# Row 30 <- former row 29 data
$ws.Range("B30").Value = 6880463
$ws.Range("E30").Value = "Erzgebirge Aue"
$ws.Range("F30").Value = "Sandhausen"
$ws.Range("G30").Value = 2
$ws.Range("H30").Value = 1
$ws.Range("I30").Value = 1
$ws.Range("J30").Value = 1
$ws.Range("K30").Value = "H"
$ws.Range("L30").Value = 2.8
$ws.Range("M30").Value = 3.4
$ws.Range("N30").Value = 2.25
$ws.Range("O30").Value = 2.8
$ws.Range("P30").Value = 3.4
$ws.Range("Q30").Value = 2.25
$ws.Range("R30").Value = 0.25
$ws.Range("S30").Value = 1.825
$ws.Range("T30").Value = 2.025
$ws.Range("U30").Value = 2.75
$ws.Range("V30").Value = 1.975
$ws.Range("W30").Value = 1.875
$ws.Range("X30").Value = 1.8
$ws.Range("Y30").Value = -1
$ws.Range("Z30").Value = -1
$ws.Range("AA30").Value = 0.825
$ws.Range("AB30").Value = -1
$ws.Range("AC30").Value = 0.4875
$ws.Range("AD30").Value = -0.5

# Row 371 <- former row 376 data
$ws.Range("B371").Value = 7280382
$ws.Range("E371").Value = "Freiburg II"
$ws.Range("F371").Value = "FC Viktoria Kln"
$ws.Range("G371").Value = 1
$ws.Range("H371").Value = 0
$ws.Range("I371").Value = 1
$ws.Range("J371").Value = 0
$ws.Range("K371").Value = "H"
$ws.Range("L371").Value = 2.8
$ws.Range("M371").Value = 3.4
$ws.Range("N371").Value = 2.375
$ws.Range("O371").Value = 2.35
$ws.Range("P371").Value = 3.6
$ws.Range("Q371").Value = 2.7
$ws.Range("R371").Value = 0
$ws.Range("S371").Value = 1.8
$ws.Range("T371").Value = 2.05
$ws.Range("U371").Value = 3
$ws.Range("V371").Value = 1.825
$ws.Range("W371").Value = 2.025
$ws.Range("X371").Value = 1.35
$ws.Range("Y371").Value = -1
$ws.Range("Z371").Value = -1
$ws.Range("AA371").Value = 0.8
$ws.Range("AB371").Value = -1
$ws.Range("AC371").Value = -1
$ws.Range("AD371").Value = 1.025

# Row 372 <- former row 371 data
$ws.Range("B372").Value = 7280802
$ws.Range("E372").Value = "SSV Ulm 1846"
$ws.Range("F372").Value = "Verl"
$ws.Range("G372").Value = 4
$ws.Range("H372").Value = 2
$ws.Range("I372").Value = 1
$ws.Range("J372").Value = 0
$ws.Range("K372").Value = "H"
$ws.Range("L372").Value = 1.75
$ws.Range("M372").Value = 3.8
$ws.Range("N372").Value = 3.9
$ws.Range("O372").Value = 1.727
$ws.Range("P372").Value = 3.9
$ws.Range("Q372").Value = 3.9
$ws.Range("R372").Value = -0.75
$ws.Range("S372").Value = 2
$ws.Range("T372").Value = 1.8
$ws.Range("U372").Value = 2.75
$ws.Range("V372").Value = 1.775
$ws.Range("W372").Value = 2.025
$ws.Range("X372").Value = 0.7270000000000001
$ws.Range("Y372").Value = -1
$ws.Range("Z372").Value = -1
$ws.Range("AA372").Value = 1
$ws.Range("AB372").Value = -1
$ws.Range("AC372").Value = 0.7749999999999999
$ws.Range("AD372").Value = -1

# Row 373 <- former row 372 data
$ws.Range("B373").Value = 7281395
$ws.Range("E373").Value = "Vfb Lubeck"
$ws.Range("F373").Value = "RotWeiss Essen"
$ws.Range("G373").Value = 3
$ws.Range("H373").Value = 3
$ws.Range("I373").Value = 2
$ws.Range("J373").Value = 2
$ws.Range("K373").Value = "D"
$ws.Range("L373").Value = 3.4
$ws.Range("M373").Value = 3.75
$ws.Range("N373").Value = 1.909
$ws.Range("O373").Value = 3.6
$ws.Range("P373").Value = 3.9
$ws.Range("Q373").Value = 1.8
$ws.Range("R373").Value = 0.75
$ws.Range("S373").Value = 1.825
$ws.Range("T373").Value = 2.025
$ws.Range("U373").Value = 3
$ws.Range("V373").Value = 1.975
$ws.Range("W373").Value = 1.875
$ws.Range("X373").Value = -1
$ws.Range("Y373").Value = 2.9
$ws.Range("Z373").Value = -1
$ws.Range("AA373").Value = 0.825
$ws.Range("AB373").Value = -1
$ws.Range("AC373").Value = 0.9750000000000001
$ws.Range("AD373").Value = -1

# Row 374 <- former row 373 data
$ws.Range("B374").Value = 7269531
$ws.Range("E374").Value = "Hallescher FC"
$ws.Range("F374").Value = "Borussia Dortmund II"
$ws.Range("G374").Value = 1
$ws.Range("H374").Value = 1
$ws.Range("I374").Value = 1
$ws.Range("J374").Value = 1
$ws.Range("K374").Value = "D"
$ws.Range("L374").Value = 2.15
$ws.Range("M374").Value = 3.75
$ws.Range("N374").Value = 2.8
$ws.Range("O374").Value = 2
$ws.Range("P374").Value = 3.9
$ws.Range("Q374").Value = 3.1
$ws.Range("R374").Value = -0.25
$ws.Range("S374").Value = 1.775
$ws.Range("T374").Value = 2.025
$ws.Range("U374").Value = 3.25
$ws.Range("V374").Value = 1.975
$ws.Range("W374").Value = 1.825
$ws.Range("X374").Value = -1
$ws.Range("Y374").Value = 2.9
$ws.Range("Z374").Value = -1
$ws.Range("AA374").Value = -0.5
$ws.Range("AB374").Value = 0.5125
$ws.Range("AC374").Value = -1
$ws.Range("AD374").Value = 0.825

# Row 375 <- former row 374 data
$ws.Range("B375").Value = 7280801
$ws.Range("E375").Value = "Erzgebirge Aue"
$ws.Range("F375").Value = "Waldhof Mannheim"
$ws.Range("G375").Value = 2
$ws.Range("H375").Value = 0
$ws.Range("I375").Value = 1
$ws.Range("J375").Value = 0
$ws.Range("K375").Value = "H"
$ws.Range("L375").Value = 2.4
$ws.Range("M375").Value = 3.5
$ws.Range("N375").Value = 2.625
$ws.Range("O375").Value = 2
$ws.Range("P375").Value = 3.8
$ws.Range("Q375").Value = 3.1
$ws.Range("R375").Value = -0.25
$ws.Range("S375").Value = 1.775
$ws.Range("T375").Value = 2.025
$ws.Range("U375").Value = 3.25
$ws.Range("V375").Value = 1.9
$ws.Range("W375").Value = 1.9
$ws.Range("X375").Value = 1
$ws.Range("Y375").Value = -1
$ws.Range("Z375").Value = -1
$ws.Range("AA375").Value = 0.7749999999999999
$ws.Range("AB375").Value = -1
$ws.Range("AC375").Value = -1
$ws.Range("AD375").Value = 0.8999999999999999

# Row 376 <- former row 375 data
$ws.Range("B376").Value = 7280383
$ws.Range("E376").Value = "Sandhausen"
$ws.Range("F376").Value = "FC Ingolstadt"
$ws.Range("G376").Value = 1
$ws.Range("H376").Value = 1
$ws.Range("I376").Value = 1
$ws.Range("J376").Value = 0
$ws.Range("K376").Value = "D"
$ws.Range("L376").Value = 2.625
$ws.Range("M376").Value = 3.5
$ws.Range("N376").Value = 2.375
$ws.Range("O376").Value = 2.55
$ws.Range("P376").Value = 3.5
$ws.Range("Q376").Value = 2.45
$ws.Range("R376").Value = 0
$ws.Range("S376").Value = 1.95
$ws.Range("T376").Value = 1.85
$ws.Range("U376").Value = 3.25
$ws.Range("V376").Value = 1.9
$ws.Range("W376").Value = 1.9
$ws.Range("X376").Value = -1
$ws.Range("Y376").Value = 2.5
$ws.Range("Z376").Value = -1
$ws.Range("AA376").Value = 0
$ws.Range("AB376").Value = 0
$ws.Range("AC376").Value = -1
$ws.Range("AD376").Value = 0.8999999999999999

# Row 379 <- former row 380 data
$ws.Range("B379").Value = 7278383
$ws.Range("E379").Value = "1860 Munich"
$ws.Range("F379").Value = "Arminia Bielefeld"
$ws.Range("G379").Value = 0
$ws.Range("H379").Value = 2
$ws.Range("I379").Value = 0
$ws.Range("J379").Value = 1
$ws.Range("K379").Value = "A"
$ws.Range("L379").Value = 2.4
$ws.Range("M379").Value = 3.5
$ws.Range("N379").Value = 2.6
$ws.Range("O379").Value = 2.625
$ws.Range("P379").Value = 3.7
$ws.Range("Q379").Value = 2.35
$ws.Range("R379").Value = 0
$ws.Range("S379").Value = 2
$ws.Range("T379").Value = 1.8
$ws.Range("U379").Value = 3
$ws.Range("V379").Value = 2
$ws.Range("W379").Value = 1.8
$ws.Range("X379").Value = -1
$ws.Range("Y379").Value = -1
$ws.Range("Z379").Value = 1.35
$ws.Range("AA379").Value = -1
$ws.Range("AB379").Value = 0.8
$ws.Range("AC379").Value = -1
$ws.Range("AD379").Value = 0.8

# Row 380 <- former row 379 data
$ws.Range("B380").Value = 7278382
$ws.Range("E380").Value = "Jahn Regensburg"
$ws.Range("F380").Value = "Saarbrucken"
$ws.Range("G380").Value = 0
$ws.Range("H380").Value = 1
$ws.Range("I380").Value = 0
$ws.Range("J380").Value = 1
$ws.Range("K380").Value = "A"
$ws.Range("L380").Value = 2.2
$ws.Range("M380").Value = 3.6
$ws.Range("N380").Value = 2.9
$ws.Range("O380").Value = 2.15
$ws.Range("P380").Value = 3.75
$ws.Range("Q380").Value = 2.875
$ws.Range("R380").Value = -0.25
$ws.Range("S380").Value = 1.975
$ws.Range("T380").Value = 1.875
$ws.Range("U380").Value = 2.75
$ws.Range("V380").Value = 1.875
$ws.Range("W380").Value = 1.975
$ws.Range("X380").Value = -1
$ws.Range("Y380").Value = -1
$ws.Range("Z380").Value = 1.875
$ws.Range("AA380").Value = -1
$ws.Range("AB380").Value = 0.875
$ws.Range("AC380").Value = -1
$ws.Range("AD380").Value = 0.9750000000000001
